$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 117, shifting existing rows 117:144 down to 118:145
$ws.Rows.Item(117).Insert()

# Populate the new row 117. Columns not explicitly changed keep the same
# values as the row that used to be at 117 (now at 118).
$ws.Range("A117").Value = 5
$ws.Range("B117").Value = "Macroferia Regional de Talca"
$ws.Range("C117").Value = "Maule"
$ws.Range("D117").Value = 45218
$ws.Range("E117").Value = 7
$ws.Range("F117").Value = 100112022
$ws.Range("G117").Value = "Arveja Verde"
$ws.Range("H117").Value = "Sin especificar"
$ws.Range("I117").Value = "Primera"
$ws.Range("J117").Value = 200
$ws.Range("K117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("M117").Value = 25000
$ws.Range("N117").Value = "$/saco 25 kilos"
$ws.Range("O117").Value = "Región del Maule"
$ws.Range("P117").Value = 1000
$ws.Range("Q117").Value = 25
$ws.Range("R117").Value = "Hortaliza"

# Match the number format used for the date column (D) on this new row to the rest of the sheet
$ws.Range("D117").NumberFormat = $ws.Range("D118").NumberFormat
